$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated CPICQ data moments (col B) for rows 82-97, per re-run of DoSMM.
$ws.Cells.Item(82, 2).Value = 1.745377659797668
$ws.Cells.Item(83, 2).Value = 1.777267575263977
$ws.Cells.Item(84, 2).Value = 1.897096157073975
$ws.Cells.Item(85, 2).Value = 2.071507215499878
$ws.Cells.Item(86, 2).Value = 2.163610219955444
$ws.Cells.Item(87, 2).Value = 2.218834638595581
$ws.Cells.Item(88, 2).Value = 2.216459274291992
$ws.Cells.Item(89, 2).Value = 2.213071584701538
$ws.Cells.Item(90, 2).Value = 2.004358291625977
$ws.Cells.Item(91, 2).Value = 1.729343056678772
$ws.Cells.Item(92, 2).Value = 1.680276989936829
$ws.Cells.Item(93, 2).Value = 1.76869010925293
$ws.Cells.Item(94, 2).Value = 2.065015554428101
$ws.Cells.Item(95, 2).Value = 2.267652988433838
$ws.Cells.Item(96, 2).Value = 2.261080026626587
$ws.Cells.Item(97, 2).Value = 2.203789710998535
